$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text must be forced to stay as text
# (matching the source inlineStr type) by setting NumberFormat to "@" first.
$textCells = @("D5","D7","D9","D10","D11","D14","D15","D17","D20","D21","D22","D23","D26","D27","D28","D29","D31","D32","D33","D35","D40","D41","D42","D44","D45","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value2 = '25.675.57'
$ws.Range('E2').Value2 = '  -1.26%  '
$ws.Range('D3').Value2 = '1.592.90'
$ws.Range('E3').Value2 = '  -2.76%  '
$ws.Range('E4').Value2 = '  +0.06%  '
$ws.Range('D5').Value2 = '208.33'
$ws.Range('E5').Value2 = '  -3.00%  '
$ws.Range('E6').Value2 = '  +0.12%  '
$ws.Range('D7').Value2 = '0.478'
$ws.Range('E7').Value2 = '  -5.13%  '
$ws.Range('E8').Value2 = '  -2.63%  '
$ws.Range('D9').Value2 = '0.0607'
$ws.Range('E9').Value2 = '  -2.73%  '
$ws.Range('D10').Value2 = '17.67'
$ws.Range('E10').Value2 = '  -4.65%  '
$ws.Range('D11').Value2 = '0.0783'
$ws.Range('E11').Value2 = '  -0.85%  '
$ws.Range('D12').Value2 = '1.818.36'
$ws.Range('E12').Value2 = '  -2.52%  '
$ws.Range('D13').Value2 = '1.598.90'
$ws.Range('E13').Value2 = '  -3.71%  '
$ws.Range('D14').Value2 = '4.02'
$ws.Range('E14').Value2 = '  -4.70%  '
$ws.Range('D15').Value2 = '0.505'
$ws.Range('E15').Value2 = '  -4.65%  '
$ws.Range('D16').Value2 = '25.688.20'
$ws.Range('E16').Value2 = '  -1.24%  '
$ws.Range('D17').Value2 = '60.18'
$ws.Range('E17').Value2 = '  -2.77%  '
$ws.Range('D18').Value2 = '0.0₃0710'
$ws.Range('E18').Value2 = '  -4.73%  '
$ws.Range('E19').Value2 = '  -0.16%  '
$ws.Range('D20').Value2 = '188.10'
$ws.Range('E20').Value2 = '  -2.10%  '
$ws.Range('D21').Value2 = '4.15'
$ws.Range('E21').Value2 = '  -2.37%  '
$ws.Range('D22').Value2 = '9.29'
$ws.Range('E22').Value2 = '  -4.76%  '
$ws.Range('D23').Value2 = '5.91'
$ws.Range('E23').Value2 = '  -3.19%  '
$ws.Range('E24').Value2 = '  +0.17%  '
$ws.Range('E25').Value2 = '  -4.25%  '
$ws.Range('D26').Value2 = '140.55'
$ws.Range('E26').Value2 = '  -2.23%  '
$ws.Range('D27').Value2 = '1.71'
$ws.Range('E27').Value2 = '  -4.89%  '
$ws.Range('B28').Value2 = 'EthereumClassic'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value2 = '14.92'
$ws.Range('E28').Value2 = '  -2.26%  '
$ws.Range('B29').Value2 = 'Cosmos'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value2 = '6.47'
$ws.Range('E29').Value2 = '  -5.52%  '
$ws.Range('E30').Value2 = '  -4.39%  '
$ws.Range('D31').Value2 = '0.0463'
$ws.Range('E31').Value2 = '  -4.35%  '
$ws.Range('D32').Value2 = '3.05'
$ws.Range('E32').Value2 = '  -3.35%  '
$ws.Range('D33').Value2 = '2.98'
$ws.Range('E33').Value2 = '  -5.16%  '
$ws.Range('E34').Value2 = '  -0.43%  '
$ws.Range('D35').Value2 = '1.46'
$ws.Range('E35').Value2 = '  -2.50%  '
$ws.Range('D36').Value2 = '1.089.87'
$ws.Range('E36').Value2 = '  -4.31%  '
$ws.Range('E37').Value2 = '  -3.54%  '
$ws.Range('E38').Value2 = '  -0.38%  '
$ws.Range('E39').Value2 = '  -2.86%  '
$ws.Range('D40').Value2 = '0.783'
$ws.Range('E40').Value2 = '  -9.85%  '
$ws.Range('D41').Value2 = '0.493'
$ws.Range('E41').Value2 = '  -5.46%  '
$ws.Range('D42').Value2 = '95.06'
$ws.Range('E42').Value2 = '  -3.50%  '
$ws.Range('D43').Value2 = '1.731.73'
$ws.Range('E43').Value2 = '  -2.46%  '
$ws.Range('D44').Value2 = '5.04'
$ws.Range('E44').Value2 = '  -3.86%  '
$ws.Range('D45').Value2 = '0.738'
$ws.Range('E45').Value2 = '  -5.23%  '
$ws.Range('D46').Value2 = '0.0₆0108'
$ws.Range('E46').Value2 = '  -6.11%  '
$ws.Range('D47').Value2 = '52.96'
$ws.Range('E47').Value2 = '  -4.22%  '
$ws.Range('D48').Value2 = '0.0510'
$ws.Range('E48').Value2 = '  -3.61%  '
$ws.Range('E49').Value2 = '  -1.20%  '
$ws.Range('E50').Value2 = '  -5.25%  '
$ws.Range('D51').Value2 = '1.01'
$ws.Range('E51').Value2 = '  +0.00%  '
